$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'267.47"
$ws.Range("D3").Value = "'22.94"
$ws.Range("D4").Value = "'6.315"
$ws.Range("D5").Value = "'0.06195"
$ws.Range("D7").Value = "'6.687"
$ws.Range("D8").Value = "'1.389"
$ws.Range("D9").Value = "'0.8324"
$ws.Range("D10").Value = "'0.01362"
$ws.Range("D11").Value = "'0.1605"
$ws.Range("D12").Value = "'0.08287"
$ws.Range("D14").Value = "'0.03165"
$ws.Range("D15").Value = "'0.09283"
$ws.Range("D16").Value = "'3.910"
$ws.Range("D17").Value = "'0.001737"
$ws.Range("D18").Value = "'0.04859"
$ws.Range("D19").Value = "'0.006246"
$ws.Range("D20").Value = "'0.005375"
$ws.Range("D21").Value = "'0.001091"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.769"
$ws.Range("D24").Value = "'2.369"
$ws.Range("D40").Value = "'0.04671"
$ws.Range("D41").Value = "'0.006900"
$ws.Range("D42").Value = "'0.1154"
$ws.Range("D43").Value = "'0.003461"
$ws.Range("D44").Value = "'0.01219"
$ws.Range("D45").Value = "'0.00006213"
$ws.Range("D47").Value = "'0.7001"
$ws.Range("D48").Value = "'0.1742"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01241"
